$wb = $excel.ActiveWorkbook

# The existing "Strommessung" sheet is duplicated; the duplicate becomes the
# new first sheet "Strommessung_Flieger" while the original is renamed to
# "Strommessung_Test".
$src = $wb.Worksheets.Item("Strommessung")
$src.Copy($src) | Out-Null

$newSheet = $wb.Worksheets.Item("Strommessung (2)")
$oldSheet = $wb.Worksheets.Item("Strommessung")

$newSheet.Name = "Strommessung_Flieger"
$oldSheet.Name = "Strommessung_Test"

# On the new "Flieger" sheet, change the Windings value from 5 to 2.
$newSheet.Range("C2").Value = 2

# Make the new sheet the active one, with C2 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("C2").Select() | Out-Null

# On the "Test" sheet (former "Strommessung"), the selection moved to C35.
$oldSheet.Range("C35").Select() | Out-Null

# "Spannungsmessung" keeps its content and selection (F8); it is simply no
# longer the active tab since "Strommessung_Flieger" is active now.
$vSheet = $wb.Worksheets.Item("Spannungsmessung")
$vSheet.Range("F8").Select() | Out-Null

$newSheet.Activate() | Out-Null
